# Add a new "2022-Q3" sheet (right after "总计") with fund holdings data,
# and add a corresponding summary row on the "总计" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update the "总计" (summary) sheet: insert a new row 2 for 2022-Q3
#    and shift the existing data rows down, renumbering the index column.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Rows.Item(2).Insert()
$summary.Rows.Item(2).ClearFormats()

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 4
$summary.Cells.Item(2, 4).Value = 1.9

# Copy the row-index cell style (bold / centered / bordered) from row 3
# onto the newly inserted row 2 so it matches the other index cells.
$summary.Cells.Item(3, 1).Copy()
$summary.Cells.Item(2, 1).PasteSpecial(-4122)

# Renumber the index column (A) for the rows that shifted down.
for ($r = 3; $r -le 9; $r++) {
    $summary.Cells.Item($r, 1).Value = $r - 2
}

# ---------------------------------------------------------------------
# 2. Insert the new "2022-Q3" worksheet right after "总计" and populate
#    it with the fund holdings detail data.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $summary)
$q3.Name = "2022-Q3"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q3.Cells.Item(1, $i + 2).Value = $headers[$i]
}

$rows = @(
    @("000011", "华夏大盘精选混合A", "40.09", "89.60", "4.62", "1.8522", 8),
    @("006973", "太平睿盈混合A",     "3.84",  "28.79", "0.72", "0.0276", 10),
    @("012628", "华夏大盘精选混合C", "0.17",  "89.60", "4.62", "0.0079", 8),
    @("007669", "太平睿盈混合C",     "1.04",  "28.79", "0.72", "0.0075", 10)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $data = $rows[$i]

    $q3.Cells.Item($r, 1).Value = $i

    # Columns B-G hold text values (fund codes must keep leading zeros, and
    # the numeric-looking figures are stored as text in this workbook).
    for ($c = 0; $c -le 4; $c++) {
        $cell = $q3.Cells.Item($r, $c + 2)
        $cell.NumberFormat = "@"
        $cell.Value = [string]$data[$c]
    }

    # Column H (仓位排名) is a real number.
    $q3.Cells.Item($r, 8).Value = $data[5]
}

# Style the index column (A) and header row (B1:H1) to match the bold,
# centered, thin-bordered look used throughout the workbook.
$styledRange = $q3.Range("A1:A5, B1:H1")
$styledRange.Font.Bold = $true
$styledRange.HorizontalAlignment = -4108
$styledRange.VerticalAlignment = -4160
$styledRange.Borders.LineStyle = 1

# ---------------------------------------------------------------------
# 3. Restore the originally active sheet (2020-Q4 was selected before).
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2020-Q4").Activate()
